$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the "last row" formatting (currently on row 26) onto row 19 before we
# delete rows, so the surviving 4th data row ends up with the bottom-border style
# that belongs to the final row of the (now shorter) table.
$ws.Range("B26:J26").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- Remove the rows that are no longer part of the table (old rows 20-26: the
# rest of the JAIME ENRIQUE MERCADO MESTRE periods, the old JONAHATAN row, and the
# old ALBERTO rows). This shifts everything below up by 7 rows, moving the
# signature/footer block from rows 31-32 to rows 24-25.
$ws.Range("A20:A26").EntireRow.Delete()

# --- Update the header summary figures.
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4
$ws.Range("E11").Value = 194454

# --- Rewrite the 4 data rows (16-19) with the new worker/period data.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1101814832"
$ws.Range("D16").Value = "ALBERTO LUIS YEPES ZABALA"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 52000
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1101814832"
$ws.Range("D17").Value = "ALBERTO LUIS YEPES ZABALA"
$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1101814832"
$ws.Range("D18").Value = "ALBERTO LUIS YEPES ZABALA"
$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 46800
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1051893418"
$ws.Range("D19").Value = "JONAHATAN LEONARDO YANCES MUÑOZ"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 43654
$ws.Range("G19").Value = 1423500
